# fix: the SOP form of const nodes
#
# Row 6 ("c1908") measurement values were updated (areas recalculated),
# which also updates the dependent average-ratio formulas in R6/R15.
# C6 also loses its old "flagged" (red-font) formatting and is restyled
# to match the plain numeric style used by the rest of the row.
# Finally, the saved selection moves to J6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C6 previously carried a distinct ("highlighted") cell style (red font).
# Re-format it to match its neighbour B6's plain style before writing the
# new value, so it converges back onto the shared/common style instead of
# keeping its own one-off style slot.
$ws.Range("B6").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C6").Value = 1615
$ws.Range("D6").Value = 1600
$ws.Range("F6").Value = 1600
$ws.Range("H6").Value = 1600
$ws.Range("J6").Value = 1579
$ws.Range("L6").Value = 1579
$ws.Range("N6").Value = 1523
$ws.Range("P6").Value = 1483

# Move/save the active selection to J6 (matches the saved sheetView).
$ws.Range("J6").Select() | Out-Null
